# Update the five divided-by practice tables (one table, 20 rows,
# content on rows 1/5/9/13/17) replacing each equation cell's text
# with the new value from the commit's regenerated data set.
# Using Cell.Range scoped Find/Replace (rather than a document-wide
# replace) because some equations such as "59÷7=" occur more than
# once and must map to different replacements depending on position.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$changes = @(
    @{ Row = 1;  Col = 1; Old = "28÷3="; New = "66÷7=" },
    @{ Row = 1;  Col = 2; Old = "99÷6="; New = "29÷3=" },
    @{ Row = 1;  Col = 3; Old = "42÷4="; New = "14÷8=" },
    @{ Row = 1;  Col = 4; Old = "33÷5="; New = "52÷2=" },
    @{ Row = 1;  Col = 5; Old = "96÷9="; New = "67÷9=" },

    @{ Row = 5;  Col = 1; Old = "49÷8="; New = "10÷9=" },
    @{ Row = 5;  Col = 2; Old = "22÷3="; New = "38÷5=" },
    @{ Row = 5;  Col = 3; Old = "89÷5="; New = "10÷2=" },
    @{ Row = 5;  Col = 4; Old = "14÷2="; New = "51÷7=" },
    @{ Row = 5;  Col = 5; Old = "91÷7="; New = "48÷3=" },

    @{ Row = 9;  Col = 1; Old = "13÷9="; New = "21÷8=" },
    @{ Row = 9;  Col = 2; Old = "55÷8="; New = "63÷9=" },
    @{ Row = 9;  Col = 3; Old = "17÷9="; New = "75÷6=" },
    @{ Row = 9;  Col = 4; Old = "27÷7="; New = "62÷4=" },
    @{ Row = 9;  Col = 5; Old = "64÷4="; New = "42÷3=" },

    @{ Row = 13; Col = 1; Old = "59÷7="; New = "26÷3=" },
    @{ Row = 13; Col = 2; Old = "93÷6="; New = "71÷6=" },
    @{ Row = 13; Col = 3; Old = "42÷9="; New = "70÷4=" },
    @{ Row = 13; Col = 4; Old = "77÷5="; New = "49÷4=" },
    @{ Row = 13; Col = 5; Old = "63÷5="; New = "44÷5=" },

    @{ Row = 17; Col = 1; Old = "48÷6="; New = "19÷7=" },
    @{ Row = 17; Col = 2; Old = "93÷5="; New = "81÷6=" },
    @{ Row = 17; Col = 3; Old = "26÷9="; New = "29÷8=" },
    @{ Row = 17; Col = 4; Old = "59÷7="; New = "35÷5=" },
    @{ Row = 17; Col = 5; Old = "67÷7="; New = "63÷7=" }
)

foreach ($chg in $changes) {
    $cellRange = $tbl.Cell($chg.Row, $chg.Col).Range
    # Replace:=1 (wdReplaceOne) -- several equations (e.g. "59÷7=") repeat
    # at multiple cells, so a "replace all" here would leak across the
    # whole document instead of staying scoped to this single cell.
    $cellRange.Find.Execute($chg.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $chg.New, 1)
}
